$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns before column I (9) to make room for the new
# taxonomy_a/taxon_a/taxonomy_b/taxon_b/taxonomy_c/taxon_c columns,
# replacing the old category_brand/category_collection/category_type/sub_category
# 4-column block with a new 6-column block.
$ws.Range("I1:J1").EntireColumn.Insert()

# Header + sample row for the new taxonomy columns I:N, written column by
# column (taxonomy_a/Brand, taxon_a/Youxi, taxonomy_b/Collection, ...) so
# new shared-string entries get interned in the same interleaved order the
# original authoring tool produced them in.
$ws.Range("I1").Value = "taxonomy_a"
$ws.Range("I2").Value = "Brand"
$ws.Range("J1").Value = "taxon_a"
$ws.Range("J2").Value = "Youxi"
$ws.Range("K1").Value = "taxonomy_b"
$ws.Range("K2").Value = "Collection"
$ws.Range("L1").Value = "taxon_b"
$ws.Range("L2").Value = "From Nature"
$ws.Range("M1").Value = "taxonomy_c"
$ws.Range("M2").Value = "Type"
$ws.Range("N1").Value = "taxon_c"
$ws.Range("N2").Value = "Home and Living"

# K2/M2 inherit the "narrow" banding style from the column insert; re-paint
# them with the same (wrap-text) format as the rest of the taxonomy row so
# the whole I2:N2 block is visually uniform, matching J2/L2/N2.
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("J2").Copy()
$ws.Range("M2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the new uniform column width used for the taxonomy block (I:N)
$ws.Range("I1:N2").ColumnWidth = 23.5

# Update the view: selection moved in the source edit
$ws.Range("U10").Select()
